$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet to reflect the new "through" date
$ws.Name = "Through 2021-09-10"

# Update the September label in column A
$ws.Range("A10").Value = "September (through 09-10)"

# Update the September row (row 10) counts
$ws.Range("B10").Value = 9
$ws.Range("C10").Value = 17
$ws.Range("D10").Value = 25
$ws.Range("E10").Value = 15
$ws.Range("F10").Value = 23
$ws.Range("G10").Value = 32
$ws.Range("H10").Value = 43

# Update the Total row (row 11) counts
$ws.Range("B11").Value = 203
$ws.Range("C11").Value = 398
$ws.Range("D11").Value = 576
$ws.Range("E11").Value = 505
$ws.Range("F11").Value = 372
$ws.Range("G11").Value = 816
$ws.Range("H11").Value = 1114
